# Add season-record columns (Wins, Losses, Ties) to the stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy the formatting of the last existing header cell (AC1)
# onto the three new header cells so they pick up the same bold/border/
# centered style used by the rest of row 1, then set their text.
$ws.Range("AC1").Copy($ws.Range("AD1:AF1"))
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-46: every team row gets the same season record.
for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 30).Value = 87   # AD
    $ws.Cells.Item($r, 31).Value = 75   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}

Write-Output "Added Wins/Losses/Ties columns (AD:AF) for rows 1-46"
